$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.669677734375
$ws.Range("B1").Value = 3.712538003921509
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 3.738809108734131
$ws.Range("E1").Value = 3.055070638656616
